$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric: force Text format first so
# Excel keeps them as strings (matches source workbook storing all
# D/E column data as text), then restore default styling afterwards.
$numericLookingCells = @('D5', 'D8', 'D9', 'D10', 'D14', 'D15', 'D16', 'D18', 'D20', 'D22', 'D23', 'D24', 'D25', 'D26', 'D29', 'D31', 'D35', 'D39', 'D42', 'D43', 'D45', 'D48', 'D50', 'D51')
foreach ($c in $numericLookingCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value2 = '27.436.61'
$ws.Range('E2').Value2 = '  +0.01%  '
$ws.Range('D3').Value2 = '1.636.90'
$ws.Range('E3').Value2 = '  -0.96%  '
$ws.Range('E4').Value2 = '  +0.01%  '
$ws.Range('D5').Value2 = '212.11'
$ws.Range('E5').Value2 = '  -0.92%  '
$ws.Range('E6').Value2 = '  +4.27%  '
$ws.Range('E7').Value2 = '  +0.03%  '
$ws.Range('D8').Value2 = '22.94'
$ws.Range('E8').Value2 = '  -4.94%  '
$ws.Range('D9').Value2 = '0.256'
$ws.Range('E9').Value2 = '  -2.48%  '
$ws.Range('D10').Value2 = '0.0609'
$ws.Range('E10').Value2 = '  -1.31%  '
$ws.Range('E11').Value2 = '  +1.06%  '
$ws.Range('D12').Value2 = '1.866.85'
$ws.Range('E12').Value2 = '  -1.09%  '
$ws.Range('D13').Value2 = '1.634.59'
$ws.Range('E13').Value2 = '  -1.29%  '
$ws.Range('D14').Value2 = '0.581'
$ws.Range('E14').Value2 = '  +2.95%  '
$ws.Range('D15').Value2 = '4.01'
$ws.Range('E15').Value2 = '  -2.62%  '
$ws.Range('D16').Value2 = '64.17'
$ws.Range('E16').Value2 = '  -2.29%  '
$ws.Range('D17').Value2 = '27.426.00'
$ws.Range('E17').Value2 = '  -0.09%  '
$ws.Range('D18').Value2 = '228.93'
$ws.Range('E18').Value2 = '  -2.95%  '
$ws.Range('D19').Value2 = '0.0₃0722'
$ws.Range('E19').Value2 = '  -0.65%  '
$ws.Range('D20').Value2 = '7.51'
$ws.Range('E20').Value2 = '  -0.31%  '
$ws.Range('E21').Value2 = '  +0.09%  '
$ws.Range('D22').Value2 = '4.30'
$ws.Range('E22').Value2 = '  -2.59%  '
$ws.Range('D23').Value2 = '9.64'
$ws.Range('E23').Value2 = '  +3.20%  '
$ws.Range('D24').Value2 = '1.96'
$ws.Range('E24').Value2 = '  -3.43%  '
$ws.Range('D25').Value2 = '149.17'
$ws.Range('E25').Value2 = '  +2.36%  '
$ws.Range('D26').Value2 = '6.97'
$ws.Range('E26').Value2 = '  -2.88%  '
$ws.Range('E27').Value2 = '  +1.36%  '
$ws.Range('E28').Value2 = '  +0.10%  '
$ws.Range('D29').Value2 = '15.52'
$ws.Range('E29').Value2 = '  -3.55%  '
$ws.Range('E30').Value2 = '  -1.01%  '
$ws.Range('D31').Value2 = '0.0487'
$ws.Range('E31').Value2 = '  -2.23%  '
$ws.Range('E32').Value2 = '  -0.71%  '
$ws.Range('E33').Value2 = '  +3.46%  '
$ws.Range('D34').Value2 = '1.408.45'
$ws.Range('E34').Value2 = '  -2.78%  '
$ws.Range('D35').Value2 = '1.59'
$ws.Range('E35').Value2 = '  +2.01%  '
$ws.Range('E36').Value2 = '  -1.55%  '
$ws.Range('E37').Value2 = '  -0.40%  '
$ws.Range('E38').Value2 = '  -4.32%  '
$ws.Range('D39').Value2 = '0.0167'
$ws.Range('E39').Value2 = '  -1.91%  '
$ws.Range('E40').Value2 = '  -1.82%  '
$ws.Range('E41').Value2 = '  +0.06%  '
$ws.Range('D42').Value2 = '0.820'
$ws.Range('E42').Value2 = '  +4.10%  '
$ws.Range('D43').Value2 = '5.48'
$ws.Range('E43').Value2 = '  +0.64%  '
$ws.Range('E44').Value2 = '  +0.35%  '
$ws.Range('D45').Value2 = '64.69'
$ws.Range('E45').Value2 = '  -2.33%  '
$ws.Range('D46').Value2 = '1.777.62'
$ws.Range('E46').Value2 = '  -0.95%  '
$ws.Range('E47').Value2 = '  -3.31%  '
$ws.Range('D48').Value2 = '85.98'
$ws.Range('E48').Value2 = '  -2.82%  '
$ws.Range('D49').Value2 = '0.0₆0107'
$ws.Range('E49').Value2 = '  +0.77%  '
$ws.Range('D50').Value2 = '0.0991'
$ws.Range('E50').Value2 = '  -2.03%  '
$ws.Range('D51').Value2 = '7.69'
$ws.Range('E51').Value2 = '  -1.23%  '

foreach ($c in $numericLookingCells) {
    $ws.Range($c).Style = "Normal"
}
